$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff that drives this update only ever touches the "Price" (D) and
# "Volume(1h)" (E) columns of the crypto table, refreshing each coin's
# latest quote. Every one of those cells is stored as TEXT in the
# workbook (prices use a dotted thousands-grouping like "29.804.15" and
# some keep a significant trailing zero like "242.77"/"1.000", so they
# must never be silently reinterpreted as numbers by Excel). When the new
# text looks like a plain number, force the cell format to Text first so
# Excel's automatic type detection doesn't convert it (and strip
# meaningful trailing zeros); then restore the "Normal" style so no
# stray number-format style lingers on the cell.
$updates = @(
    @{ Ref = "D2"; Value = "29.804.15" }
    @{ Ref = "E2"; Value = "  -0.42%  " }
    @{ Ref = "D3"; Value = "1.892.54" }
    @{ Ref = "E3"; Value = "  +0.05%  " }
    @{ Ref = "E4"; Value = "  -0.14%  " }
    @{ Ref = "D5"; Value = "0.8004" }
    @{ Ref = "E5"; Value = "  -3.94%  " }
    @{ Ref = "D6"; Value = "242.77" }
    @{ Ref = "E6"; Value = "  +0.45%  " }
    @{ Ref = "E7"; Value = "  -0.12%  " }
    @{ Ref = "D8"; Value = "0.3173" }
    @{ Ref = "E8"; Value = "  -2.53%  " }
    @{ Ref = "D9"; Value = "25.56" }
    @{ Ref = "E9"; Value = "  -4.16%  " }
    @{ Ref = "E10"; Value = "  +0.27%  " }
    @{ Ref = "D11"; Value = "0.08046" }
    @{ Ref = "E11"; Value = "  +0.31%  " }
    @{ Ref = "D12"; Value = "0.7678" }
    @{ Ref = "E12"; Value = "  +2.62%  " }
    @{ Ref = "D13"; Value = "1.909.97" }
    @{ Ref = "E13"; Value = "  +0.90%  " }
    @{ Ref = "D14"; Value = "5.294" }
    @{ Ref = "D15"; Value = "92.25" }
    @{ Ref = "E15"; Value = "  -0.04%  " }
    @{ Ref = "D16"; Value = "29.817.27" }
    @{ Ref = "E17"; Value = "  -1.26%  " }
    @{ Ref = "D18"; Value = "5.936" }
    @{ Ref = "E18"; Value = "  +0.36%  " }
    @{ Ref = "D19"; Value = "243.98" }
    @{ Ref = "E19"; Value = "  +0.13%  " }
    @{ Ref = "D20"; Value = "0.000007712" }
    @{ Ref = "E20"; Value = "  -0.68%  " }
    @{ Ref = "D21"; Value = "8.213" }
    @{ Ref = "E21"; Value = "  +18.70%  " }
    @{ Ref = "D22"; Value = "0.9986" }
    @{ Ref = "E22"; Value = "  -0.15%  " }
    @{ Ref = "D23"; Value = "2.140.67" }
    @{ Ref = "E23"; Value = "  -0.45%  " }
    @{ Ref = "D24"; Value = "0.9988" }
    @{ Ref = "E24"; Value = "  -0.14%  " }
    @{ Ref = "D25"; Value = "0.1675" }
    @{ Ref = "E25"; Value = "  +3.92%  " }
    @{ Ref = "D26"; Value = "9.323" }
    @{ Ref = "E26"; Value = "  +1.42%  " }
    @{ Ref = "D27"; Value = "165.73" }
    @{ Ref = "E27"; Value = "  -1.26%  " }
    @{ Ref = "E28"; Value = "  -0.87%  " }
    @{ Ref = "D29"; Value = "2.059" }
    @{ Ref = "E29"; Value = "  -1.01%  " }
    @{ Ref = "D30"; Value = "1.395" }
    @{ Ref = "E30"; Value = "  +1.75%  " }
    @{ Ref = "E31"; Value = "  +1.10%  " }
    @{ Ref = "E32"; Value = "  +3.60%  " }
    @{ Ref = "D33"; Value = "0.05657" }
    @{ Ref = "E33"; Value = "  +0.63%  " }
    @{ Ref = "D34"; Value = "4.053" }
    @{ Ref = "E34"; Value = "  -0.38%  " }
    @{ Ref = "D35"; Value = "1.262" }
    @{ Ref = "E35"; Value = "  -1.24%  " }
    @{ Ref = "D36"; Value = "0.7393" }
    @{ Ref = "E36"; Value = "  +0.90%  " }
    @{ Ref = "D37"; Value = "1.002" }
    @{ Ref = "E37"; Value = "  +0.25%  " }
    @{ Ref = "D38"; Value = "2.633" }
    @{ Ref = "E38"; Value = "  -3.12%  " }
    @{ Ref = "E39"; Value = "  -0.12%  " }
    @{ Ref = "D40"; Value = "2.779" }
    @{ Ref = "E40"; Value = "  +0.10%  " }
    @{ Ref = "E41"; Value = "  -0.05%  " }
    @{ Ref = "D42"; Value = "72.70" }
    @{ Ref = "E42"; Value = "  +1.04%  " }
    @{ Ref = "D43"; Value = "5.814" }
    @{ Ref = "E43"; Value = "  -2.12%  " }
    @{ Ref = "D44"; Value = "0.8452" }
    @{ Ref = "E44"; Value = "  +0.29%  " }
    @{ Ref = "D45"; Value = "0.9983" }
    @{ Ref = "E45"; Value = "  -0.14%  " }
    @{ Ref = "D46"; Value = "1.032.24" }
    @{ Ref = "E46"; Value = "  +3.91%  " }
    @{ Ref = "D47"; Value = "102.64" }
    @{ Ref = "E47"; Value = "  +1.35%  " }
    @{ Ref = "D48"; Value = "1.871" }
    @{ Ref = "E48"; Value = "  -0.99%  " }
    @{ Ref = "D49"; Value = "9.934" }
    @{ Ref = "E49"; Value = "  +2.65%  " }
    @{ Ref = "D50"; Value = "7.430" }
    @{ Ref = "E50"; Value = "  -2.14%  " }
    @{ Ref = "D51"; Value = "2.040.59" }
    @{ Ref = "E51"; Value = "  -0.33%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $val = $u.Value
    $looksNumeric = $val -match '^[+-]?\d+(\.\d+)?$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
